$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a range as plain TEXT (avoiding Excel's automatic
# number/date auto-conversion for numeric- or date-looking strings) while
# preserving the destination cell's existing style. We stage the text in a
# scratch cell formatted as Text, copy/paste-special (values) it into the
# real destination (this keeps the destination's own style untouched), then
# remove the scratch column entirely.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Cells.Item(1, 10)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.EntireColumn.Delete()
}

# Créditos-aula: 4 -> 2
Set-TextValue $ws.Range("B5") "2"
Set-TextValue $ws.Range("C5") "2"

# Carga horária: 60 h -> 30 h
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# Ativação: 01/01/2012 -> 01/01/2022
Set-TextValue $ws.Range("B8") "01/01/2022"
Set-TextValue $ws.Range("C8") "01/01/2022"

# Objetivos: collapse internal line breaks into a single line
$objetivos = "Introduzir ao aluno a teoria de propriedades elétricas, térmicas, magnéticas e óticas de materiais sólidos, levando emconta o aspecto microscópico da estrutura do material. Dá-se ênfase à aplicação do material de acordo com aspropriedades que ele apresenta."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Programa resumido: collapse internal line breaks into a single line
$programaResumido = "PROPRIEDADES ELETRÔNICAS: Condutividade elétrica em metais, semicondutores e isolantes.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão, Condutividade e Tensões Térmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Lasers. Aplicações."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Programa: collapse internal line breaks into a single line
$programa = "PROPRIEDADES ELETRÔNICAS:Teoria do Elétron Livre em Metais. Níveis de Energia em Sólidos. Condutividade.Supercondutividade. Semicondutividade. Isolantes (Dielétricos). Aplicações.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão Térmica. Condutividade Térmica. TensõesTérmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Fotocondutividade. Luminescência. Lasers.Fibra Ótica. Danos por Radiação. Aplicações."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa
